$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 10
$ws.Range("I8").Value = 6.6666665
$ws.Range("K8").Value = 19.9999995
$ws.Range("M8").Value = 119.0000005
$ws.Range("H100").Value = 4406.8
$ws.Range("I100").Value = 3093.3333
$ws.Range("K100").Value = 3093.3333
$ws.Range("M100").Value = -2552.3333
$ws.Range("H132").Value = 31145.412
$ws.Range("I132").Value = 1831.4667
$ws.Range("K132").Value = 5494.4001
$ws.Range("M132").Value = -2964.4001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 12755.333
$ws.Range("I2").Value = 15097.454
$ws.Range("K2").Value = 15097.454
$ws.Range("M2").Value = -14984.454
$ws.Range("H102").Value = 23812040
$ws.Range("I102").Value = 2433.4167
$ws.Range("K102").Value = 2433.4167
$ws.Range("M102").Value = -811.4167000000002
$ws.Range("H110").Value = 1884.0526
$ws.Range("I110").Value = 2133.4
$ws.Range("K110").Value = 2133.4
$ws.Range("M110").Value = -88.40000000000009
$ws.Range("H116").Value = 12755.333
$ws.Range("I116").Value = 15097.454
$ws.Range("K116").Value = 15097.454
$ws.Range("M116").Value = -12803.454
$ws.Range("H122").Value = 4612.3516
$ws.Range("I122").Value = 3687.4707
$ws.Range("J122").Value = 5398.5
$ws.Range("K122").Value = 11062.4121
$ws.Range("L122").Value = 16195.5
$ws.Range("M122").Value = -8612.4121
$ws.Range("N122").Value = -21095.5
$ws.Range("H132").Value = 3599.5715
$ws.Range("I132").Value = 2912.25
$ws.Range("J132").Value = 5317.875
$ws.Range("K132").Value = 8736.75
$ws.Range("L132").Value = 15953.625
$ws.Range("M132").Value = -6206.75
$ws.Range("N132").Value = -21013.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 12755.333
$ws.Range("I3").Value = 15097.454
$ws.Range("K3").Value = 15097.454
$ws.Range("M3").Value = -14983.454
$ws.Range("H63").Value = 75000
$ws.Range("J63").Value = 75000
$ws.Range("L63").Value = 75000
$ws.Range("N63").Value = -76372
$ws.Range("H66").Value = 75000
$ws.Range("J66").Value = 75000
$ws.Range("L66").Value = 225000
$ws.Range("N66").Value = -231864
$ws.Range("H99").Value = 135736.67
$ws.Range("I99").Value = 201605
$ws.Range("K99").Value = 201605
$ws.Range("M99").Value = -200107
$ws.Range("H134").Value = 5802.8237
$ws.Range("I134").Value = 4319.143
$ws.Range("J134").Value = 6841.4
$ws.Range("K134").Value = 12957.429
$ws.Range("L134").Value = 20524.2
$ws.Range("M134").Value = -10422.429
$ws.Range("N134").Value = -25594.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 289.52942
$ws.Range("I7").Value = 63.75
$ws.Range("J7").Value = 359
$ws.Range("K7").Value = 63.75
$ws.Range("L7").Value = 359
$ws.Range("M7").Value = 49.25
$ws.Range("N7").Value = -585
$ws.Range("H99").Value = 10019.454
$ws.Range("I99").Value = 15000
$ws.Range("J99").Value = 9521.4
$ws.Range("K99").Value = 15000
$ws.Range("L99").Value = 9521.4
$ws.Range("M99").Value = -13502
$ws.Range("N99").Value = -12517.4
$ws.Range("H126").Value = 10019.454
$ws.Range("I126").Value = 15000
$ws.Range("J126").Value = 9521.4
$ws.Range("K126").Value = 45000
$ws.Range("L126").Value = 28564.2
$ws.Range("M126").Value = -42530
$ws.Range("N126").Value = -33504.2
$ws.Range("H132").Value = 2895.5293
$ws.Range("I132").Value = 2556.5386
$ws.Range("K132").Value = 7669.6158
$ws.Range("M132").Value = -5139.6158
$ws.Range("H134").Value = 3828.7346
$ws.Range("I134").Value = 2672.6365
$ws.Range("J134").Value = 6213.1875
$ws.Range("K134").Value = 8017.9095
$ws.Range("L134").Value = 18639.5625
$ws.Range("M134").Value = -5482.9095
$ws.Range("N134").Value = -23709.5625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 27779310
$ws.Range("I140").Value = 27779310
$ws.Range("K140").Value = 83337930
$ws.Range("M140").Value = -83332750

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H32").Value = 40290
$ws.Range("J32").Value = 40290
$ws.Range("L32").Value = 40290
$ws.Range("N32").Value = -40882
$ws.Range("H70").Value = 164657.86
$ws.Range("I70").Value = 280651.5
$ws.Range("J70").Value = 9999.666999999999
$ws.Range("K70").Value = 280651.5
$ws.Range("L70").Value = 9999.666999999999
$ws.Range("M70").Value = -280381.5
$ws.Range("N70").Value = -10539.667
$ws.Range("H73").Value = 164657.86
$ws.Range("I73").Value = 280651.5
$ws.Range("J73").Value = 9999.666999999999
$ws.Range("K73").Value = 280651.5
$ws.Range("L73").Value = 9999.666999999999
$ws.Range("M73").Value = -279715.5
$ws.Range("N73").Value = -11871.667
$ws.Range("H102").Value = 2557.4546
$ws.Range("I102").Value = 1808.0526
$ws.Range("K102").Value = 1808.0526
$ws.Range("M102").Value = -186.0526
$ws.Range("H113").Value = 6428.9375
$ws.Range("I113").Value = 3995
$ws.Range("J113").Value = 10485.5
$ws.Range("K113").Value = 3995
$ws.Range("L113").Value = 10485.5
$ws.Range("M113").Value = -1825
$ws.Range("N113").Value = -14825.5
$ws.Range("H122").Value = 4769.5654
$ws.Range("I122").Value = 3873.9092
$ws.Range("J122").Value = 5590.5835
$ws.Range("K122").Value = 11621.7276
$ws.Range("L122").Value = 16771.7505
$ws.Range("M122").Value = -9171.7276
$ws.Range("N122").Value = -21671.7505
$ws.Range("H132").Value = 3488.7144
$ws.Range("I132").Value = 2783
$ws.Range("J132").Value = 4429.6665
$ws.Range("K132").Value = 8349
$ws.Range("L132").Value = 13288.9995
$ws.Range("M132").Value = -5819
$ws.Range("N132").Value = -18348.9995

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6790.316
$ws.Range("I7").Value = 8406.777
$ws.Range("K7").Value = 8406.777
$ws.Range("M7").Value = -8294.777
$ws.Range("H93").Value = 296580.56
$ws.Range("I93").Value = 2610.15
$ws.Range("J93").Value = 716538.3
$ws.Range("K93").Value = 2610.15
$ws.Range("L93").Value = 716538.3
$ws.Range("M93").Value = -1362.15
$ws.Range("N93").Value = -719034.3
$ws.Range("H126").Value = 6790.316
$ws.Range("I126").Value = 8406.777
$ws.Range("K126").Value = 25220.331
$ws.Range("M126").Value = -22750.331
$ws.Range("H132").Value = 4354.6924
$ws.Range("I132").Value = 3522.1875
$ws.Range("J132").Value = 5686.7
$ws.Range("K132").Value = 10566.5625
$ws.Range("L132").Value = 17060.1
$ws.Range("M132").Value = -8036.5625
$ws.Range("N132").Value = -22120.1

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").ClearContents()
$ws.Range("I81").Value = 10166.546
$ws.Range("K81").Value = 20333.092
$ws.Range("M81").Value = -19272.092
$ws.Range("I84").Value = 10166.546
$ws.Range("K84").Value = 101665.46
$ws.Range("M84").Value = -96361.46000000001
$ws.Range("H109").Value = 28826.357
$ws.Range("J109").Value = 28826.357
$ws.Range("L109").Value = 28826.357
$ws.Range("N109").Value = -31600.357
$ws.Range("H132").Value = 2569.0256
$ws.Range("I132").Value = 1990.8182
$ws.Range("K132").Value = 5972.4546
$ws.Range("M132").Value = -3442.4546
